$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.183.92"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "'1.677.40"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'211.29"
$ws.Range("E5").Value = "  -3.73%  "
$ws.Range("D6").Value = "'0.5243"
$ws.Range("E6").Value = "  -4.88%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").Value = "'0.2657"
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").Value = "'0.06295"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").Value = "'21.42"
$ws.Range("D11").Value = "'0.07545"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "'1.673.88"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'4.455"
$ws.Range("D14").Value = "'0.5644"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'0.000008033"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").Value = "'66.62"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "'26.256.88"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D19").Value = "'4.835"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'188.03"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("E21").Value = "  -5.10%  "
$ws.Range("D22").Value = "'6.207"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'149.40"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'0.1247"
$ws.Range("E25").Value = "  -6.13%  "
$ws.Range("D26").Value = "'7.593"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").Value = "'16.01"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'0.06177"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'1.360"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "'1.283"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").Value = "'3.495"
$ws.Range("E32").Value = "  -4.83%  "
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "'0.6065"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'2.744"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'6.098"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "'0.01613"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").Value = "'1.082.05"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").Value = "'0.8697"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "'100.07"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "'1.828.75"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "'56.31"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "'0.9989"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "'7.980"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "'0.05240"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'0.4258"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  -2.45%  "
